# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet right before the "总计" (totals) sheet,
#    with the same per-fund holdings layout used by the other quarter
#    sheets (2021-Q2 / 2021-Q3 / 2021-Q4).
# 2. Recreate the "总计" sheet with a new leading row summarising the
#    2022-Q1 quarter, pushing the previous rows down by one.

$wb = $excel.ActiveWorkbook

# --- locate existing sheets -------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$totalOld = $wb.Worksheets.Item("总计")

# Grab the formats we want to reuse (header row style + index-column style)
# from the "2021-Q4" sheet, which already carries style id 2 (bold, bordered,
# centered) on B1:H1 and on its A column.
$fmtHeaderSrc = $q4.Range("B1:H1")
$fmtIndexSrc = $q4.Range("A2")

# Remove the old "总计" sheet -- it will be rebuilt from scratch further down
# so its data rows line up with fresh styling, same as the other sheets.
$totalOld.Delete()

# --- 1. create the "2022-Q1" sheet, right after "2021-Q4" ------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# holdings rows (fund code / name / size / stock position / position pct /
# held value / position rank)
$q1Rows = @(
    @("501080", "中金科创主题 3 年封闭运作灵活配置混合", "16.00", "79.20", "1.79", "0.2864", 5),
    @("159855", "银华中证影视主题ETF", "0.96", "97.27", "8.03", "0.0771", 2),
    @("159804", "国寿安保国证创业板中盘精选88ETF", "2.10", "98.79", "2.03", "0.0426", 7),
    @("516620", "国泰中证影视主题ETF", "0.33", "96.08", "8.00", "0.0264", 2)
)

$r = 2
foreach ($row in $q1Rows) {
    $q1.Range("A$r").Value = ($r - 2)
    $q1.Range("B$r").NumberFormat = "@"
    $q1.Range("B$r").Value = $row[0]
    $q1.Range("C$r").Value = $row[1]
    $q1.Range("D$r").NumberFormat = "@"
    $q1.Range("D$r").Value = $row[2]
    $q1.Range("E$r").NumberFormat = "@"
    $q1.Range("E$r").Value = $row[3]
    $q1.Range("F$r").NumberFormat = "@"
    $q1.Range("F$r").Value = $row[4]
    $q1.Range("G$r").NumberFormat = "@"
    $q1.Range("G$r").Value = $row[5]
    $q1.Range("H$r").Value = $row[6]
    $r = $r + 1
}

# re-apply the shared header / index-column formatting so the new sheet
# matches the look of the other quarter sheets
$fmtHeaderSrc.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$fmtIndexSrc.Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)

# --- 2. rebuild the "总计" sheet after "2022-Q1" ----------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 4, 0.43),
    @("2021-Q4", 3, 0.43),
    @("2021-Q3", 1, 0.09),
    @("2021-Q2", 1, 0.1)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Range("A$r").Value = ($r - 2)
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
    $r = $r + 1
}

$fmtHeaderSrc2 = $q4.Range("B1:D1")
$fmtHeaderSrc2.Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$fmtIndexSrc.Copy()
$total.Range("A2:A5").PasteSpecial(-4122)

Write-Output "done"
